$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 8 -> 421f3ead-d769-4ac9-9727-95e04562c96a.md : Ready for handoff -> In Translation
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$zhcn.Range("C8").Value = "In Translation"
$dede.Range("C8").Value = "In Translation"

# Row 9 -> 625e8b32-19e5-4ca9-9d71-5bad46a05340.md : Ready for handoff -> In Translation
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"
$zhcn.Range("C9").Value = "In Translation"
$dede.Range("C9").Value = "In Translation"
